# Publish terminology IG 2.0.0
# - bump Version to 1.1.0
# - Title value now matches the Name value (MedComObservationResultGroup)
# - fix typo in Description: "intende" -> "intended"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B3").Value = "1.1.0"
$ws.Range("B5").Value = $ws.Range("B4").Value2
$ws.Range("B12").Value = "The observation result group are intended to be used to sort the individual analyzes."
